$d = $word.ActiveDocument

# Locate the last bullet item in the "Out-of-Scope" list of the Scope and
# Exclusion table ("App for specific vendors.") so we can add a new list
# item, "Financial indicators.", right after it (same bulleted/numbered
# list, numId 7).
$rng = $d.Content
$found = $rng.Find.Execute("App for specific vendors.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph 'App for specific vendors.'"
}

# $rng now spans exactly the matched text; collapse to its end point so we
# can insert immediately after it (right before the paragraph mark).
$insertionPoint = $d.Range($rng.End, $rng.End)

# Build a new list-item paragraph that mirrors the formatting of the
# existing list items in this list (widowControl off, numPr ilvl 0 /
# numId 7, no paragraph borders, gray 24pt text).
$newParagraphXml = '<w:p><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:rPr><w:color w:val="434343"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="434343"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Financial indicators.</w:t></w:r></w:p>'

[void]$insertionPoint.InsertXML($newParagraphXml)
